# Rename the worksheet from "Sheet1" to "Quiz Template"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Quiz Template"

# The single shared "Incorrect answer" header used for columns H, I and J
# is replaced with three distinct, numbered headers so each incorrect
# answer column is labeled individually.
$ws.Range("H2").Value = "Incorrect answer 1"
$ws.Range("I2").Value = "Incorrect answer 2"
$ws.Range("J2").Value = "Incorrect answer 3"
